# Auto-generated Excel COM-interop edit script
# Applies scheduled market-data refresh updates to multiple sheets (columns H-N)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 286.08
$ws.Range("I15").Value = 286.08
$ws.Range("K15").Value = 858.24
$ws.Range("M15").Value = -689.24

$ws.Range("H44").Value = 95000
$ws.Range("J44").Value = 95000
$ws.Range("L44").Value = 95000
$ws.Range("N44").Value = -95924

$ws.Range("H111").Value = 1162.8182
$ws.Range("I111").Value = 1114.75
$ws.Range("J111").Value = 1190.2858
$ws.Range("K111").Value = 3344.25
$ws.Range("L111").Value = 3570.8574
$ws.Range("M111").Value = -277.25
$ws.Range("N111").Value = -9704.857400000001

$ws.Range("H135").Value = 3904.037
$ws.Range("I135").Value = 1753.619
$ws.Range("J135").Value = 11430.5
$ws.Range("K135").Value = 15782.571
$ws.Range("L135").Value = 102874.5
$ws.Range("M135").Value = -13247.571
$ws.Range("N135").Value = -107944.5

$ws.Range("H137").Value = 1050
$ws.Range("I137").Value = 775
$ws.Range("J137").Value = 1737.5
$ws.Range("K137").Value = 2325
$ws.Range("L137").Value = 5212.5
$ws.Range("M137").Value = 225
$ws.Range("N137").Value = -10312.5

$ws.Range("H141").Value = 9345.375
$ws.Range("I141").Value = 14475
$ws.Range("J141").Value = 4215.75
$ws.Range("K141").Value = 43425
$ws.Range("L141").Value = 12647.25
$ws.Range("M141").Value = -38245
$ws.Range("N141").Value = -23007.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5788.844
$ws.Range("I31").Value = 4479.086
$ws.Range("J31").Value = 10373
$ws.Range("K31").Value = 4479.086
$ws.Range("L31").Value = 10373
$ws.Range("M31").Value = -4184.086
$ws.Range("N31").Value = -10963

$ws.Range("H34").Value = 5788.844
$ws.Range("I34").Value = 4479.086
$ws.Range("J34").Value = 10373
$ws.Range("K34").Value = 4479.086
$ws.Range("L34").Value = 10373
$ws.Range("M34").Value = -4277.086
$ws.Range("N34").Value = -10777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 529.7273
$ws.Range("J4").Value = 820
$ws.Range("L4").Value = 2460
$ws.Range("N4").Value = -2684

$ws.Range("H5").Value = 904.1667
$ws.Range("I5").Value = 631.6667
$ws.Range("J5").Value = 2266.6667
$ws.Range("K5").Value = 1895.0001
$ws.Range("L5").Value = 6800.000100000001
$ws.Range("M5").Value = -1783.0001
$ws.Range("N5").Value = -7024.000100000001

$ws.Range("H9").Value = 2249.875
$ws.Range("J9").Value = 2249.875
$ws.Range("L9").Value = 6749.625
$ws.Range("N9").Value = -7197.625

$ws.Range("H15").Value = 518.5714
$ws.Range("I15").Value = 10
$ws.Range("J15").Value = 603.3333
$ws.Range("K15").Value = 30
$ws.Range("L15").Value = 1809.9999
$ws.Range("M15").Value = 110
$ws.Range("N15").Value = -2089.9999

$ws.Range("H20").Value = 1862.5
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 1862.5
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 5587.5
$ws.Range("N20").Value = -6041.5
$ws.Range("M20").ClearContents()

$ws.Range("H21").Value = 3272
$ws.Range("I21").Value = 500
$ws.Range("J21").Value = 3618.5
$ws.Range("K21").Value = 1500
$ws.Range("L21").Value = 10855.5
$ws.Range("M21").Value = -1327
$ws.Range("N21").Value = -11201.5

$ws.Range("H26").Value = 1088
$ws.Range("I26").Value = 20
$ws.Range("K26").Value = 60
$ws.Range("M26").Value = 228

$ws.Range("H40").Value = 185.05556
$ws.Range("I40").Value = 157.95833
$ws.Range("J40").Value = 239.25
$ws.Range("K40").Value = 631.83332
$ws.Range("L40").Value = 957
$ws.Range("M40").Value = -562.83332
$ws.Range("N40").Value = -1095

$ws.Range("H46").Value = 26105.334
$ws.Range("I46").Value = 655.1429000000001
$ws.Range("J46").Value = 48374.25
$ws.Range("K46").Value = 1965.4287
$ws.Range("L46").Value = 145122.75
$ws.Range("M46").Value = -1874.4287
$ws.Range("N46").Value = -145304.75

$ws.Range("H51").Value = 1739.8
$ws.Range("I51").Value = 233
$ws.Range("J51").Value = 4000
$ws.Range("K51").Value = 699
$ws.Range("L51").Value = 12000
$ws.Range("M51").Value = -239
$ws.Range("N51").Value = -12920

$ws.Range("H57").Value = 45460600
$ws.Range("I57").Value = 500003260
$ws.Range("J57").Value = 6335
$ws.Range("K57").Value = 1500009780
$ws.Range("L57").Value = 19005
$ws.Range("M57").Value = -1500009221
$ws.Range("N57").Value = -20123

$ws.Range("H58").Value = 1847.5
$ws.Range("I58").Value = 350
$ws.Range("J58").Value = 2346.6667
$ws.Range("K58").Value = 1050
$ws.Range("L58").Value = 7040.000100000001
$ws.Range("M58").Value = -922
$ws.Range("N58").Value = -7296.000100000001

$ws.Range("H135").Value = 904.1667
$ws.Range("I135").Value = 631.6667
$ws.Range("J135").Value = 2266.6667
$ws.Range("K135").Value = 5685.0003
$ws.Range("L135").Value = 20400.0003
$ws.Range("M135").Value = -3150.0003
$ws.Range("N135").Value = -25470.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 38462784
$ws.Range("I113").Value = 1053.7858
$ws.Range("J113").Value = 83334800
$ws.Range("K113").Value = 1053.7858
$ws.Range("L113").Value = 83334800
$ws.Range("M113").Value = 1116.2142
$ws.Range("N113").Value = -83339140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2246.85
$ws.Range("I61").Value = 1472.8462
$ws.Range("J61").Value = 3684.2856
$ws.Range("K61").Value = 1472.8462
$ws.Range("L61").Value = 3684.2856
$ws.Range("M61").Value = -1270.8462
$ws.Range("N61").Value = -4088.2856

$ws.Range("H100").Value = 3082.7058
$ws.Range("I100").Value = 2176.875
$ws.Range("K100").Value = 2176.875
$ws.Range("M100").Value = -1635.875

$ws.Range("H113").Value = 2246.85
$ws.Range("I113").Value = 1472.8462
$ws.Range("J113").Value = 3684.2856
$ws.Range("K113").Value = 1472.8462
$ws.Range("L113").Value = 3684.2856
$ws.Range("M113").Value = 697.1538
$ws.Range("N113").Value = -8024.2856

$ws.Range("H132").Value = 5169.5835
$ws.Range("I132").Value = 5510.4736
$ws.Range("J132").Value = 3874.2
$ws.Range("K132").Value = 16531.4208
$ws.Range("L132").Value = 11622.6
$ws.Range("M132").Value = -14001.4208
$ws.Range("N132").Value = -16682.6

$ws.Range("H136").Value = 6268
$ws.Range("I136").Value = 1902.4
$ws.Range("J136").Value = 11725
$ws.Range("K136").Value = 5707.200000000001
$ws.Range("L136").Value = 35175
$ws.Range("M136").Value = -3157.200000000001
$ws.Range("N136").Value = -40275
